$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: cell reference -> new text value
# Values must remain text (matching original inlineStr formatting),
# so we force Text number format before assignment and restore the
# default "Normal" style afterward to avoid altering cell styling.
$updates = @(
    @{Cell="D2"; Value="301.16"},
    @{Cell="E2"; Value="-4.44%"},
    @{Cell="D3"; Value="35.13"},
    @{Cell="E3"; Value="-0.89%"},
    @{Cell="D4"; Value="5.042"},
    @{Cell="E4"; Value="-1.86%"},
    @{Cell="D5"; Value="0.07964"},
    @{Cell="E5"; Value="-2.03%"},
    @{Cell="D6"; Value="1.892"},
    @{Cell="E6"; Value="-11.39%"},
    @{Cell="D7"; Value="7.797"},
    @{Cell="E7"; Value="-2.06%"},
    @{Cell="D8"; Value="4.052"},
    @{Cell="E8"; Value="-2.32%"},
    @{Cell="D9"; Value="2.914"},
    @{Cell="E9"; Value="5.50%"},
    @{Cell="D10"; Value="0.9231"},
    @{Cell="E10"; Value="-0.77%"},
    @{Cell="D11"; Value="0.1252"},
    @{Cell="E11"; Value="24.30%"},
    @{Cell="D12"; Value="0.1855"},
    @{Cell="E12"; Value="-1.03%"},
    @{Cell="D13"; Value="0.09968"},
    @{Cell="E13"; Value="9.42%"},
    @{Cell="D14"; Value="0.03582"},
    @{Cell="E14"; Value="-0.60%"},
    @{Cell="D15"; Value="0.09849"},
    @{Cell="E15"; Value="-0.62%"},
    @{Cell="D16"; Value="0.001396"},
    @{Cell="E16"; Value="-2.90%"},
    @{Cell="D17"; Value="0.005918"},
    @{Cell="E17"; Value="4.14%"},
    @{Cell="D18"; Value="3.506"},
    @{Cell="E18"; Value="1.17%"},
    @{Cell="D19"; Value="0.3399"},
    @{Cell="E19"; Value="-0.37%"},
    @{Cell="E20"; Value="-2.90%"},
    @{Cell="D21"; Value="5.044"},
    @{Cell="E21"; Value="-1.03%"},
    @{Cell="D23"; Value="0.04498"},
    @{Cell="E23"; Value="-1.35%"},
    @{Cell="E24"; Value="-2.64%"},
    @{Cell="D25"; Value="0.004783"},
    @{Cell="E25"; Value="1.65%"},
    @{Cell="D26"; Value="0.0001252"},
    @{Cell="E26"; Value="0.05%"},
    @{Cell="D27"; Value="0.0003003"},
    @{Cell="E27"; Value="-33.32%"},
    @{Cell="D39"; Value="0.01882"},
    @{Cell="E39"; Value="-3.96%"},
    @{Cell="D40"; Value="0.04718"},
    @{Cell="E40"; Value="-2.91%"},
    @{Cell="D41"; Value="0.007515"},
    @{Cell="E41"; Value="-2.55%"},
    @{Cell="D42"; Value="0.01025"},
    @{Cell="E42"; Value="30.73%"},
    @{Cell="E43"; Value="-4.83%"},
    @{Cell="D44"; Value="0.002113"},
    @{Cell="E44"; Value="-2.71%"},
    @{Cell="D45"; Value="0.01068"},
    @{Cell="E45"; Value="-9.65%"},
    @{Cell="D46"; Value="0.00006255"},
    @{Cell="E46"; Value="-5.38%"},
    @{Cell="D47"; Value="0.00000000751"},
    @{Cell="E47"; Value="-0.01%"},
    @{Cell="E48"; Value="70.97%"},
    @{Cell="D50"; Value="0.00002102"},
    @{Cell="E50"; Value="-0.01%"},
    @{Cell="D51"; Value="0.0002002"},
    @{Cell="E51"; Value="-0.01%"}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
